$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Company Name"
$ws.Range("C1").Value = "Invoice Number"
$ws.Range("D1").Value = "Total Amount"
$ws.Range("B1").Value = "Customer ID"
$ws.Range("E1").Value = "Date"

$ws.Range("E3").Select()
